$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H2").Value = 655.4
$ws.Range("I2").Value = 655.4
$ws.Range("K2").Value = 655.4
$ws.Range("M2").Value = -542.4

$ws.Range("H18").Value = 368.42856
$ws.Range("I18").Value = 368.42856
$ws.Range("K18").Value = 368.42856
$ws.Range("M18").Value = -84.42856

$ws.Range("H76").Value = 4417.1113
$ws.Range("I76").Value = 3388.889
$ws.Range("J76").Value = 5445.3335
$ws.Range("K76").Value = 3388.889
$ws.Range("L76").Value = 5445.3335
$ws.Range("M76").Value = -3073.889
$ws.Range("N76").Value = -6075.3335

$ws.Range("H79").Value = 4417.1113
$ws.Range("I79").Value = 3388.889
$ws.Range("J79").Value = 5445.3335
$ws.Range("K79").Value = 3388.889
$ws.Range("L79").Value = 5445.3335
$ws.Range("M79").Value = -2296.889
$ws.Range("N79").Value = -7629.3335

$ws.Range("H116").Value = 6670.7144
$ws.Range("I116").Value = 12756.111
$ws.Range("J116").Value = 2106.6667
$ws.Range("K116").Value = 12756.111
$ws.Range("L116").Value = 2106.6667
$ws.Range("M116").Value = -9314.111000000001
$ws.Range("N116").Value = -8990.6667

$ws.Range("H121").Value = 1471
$ws.Range("I121").Value = 845
$ws.Range("J121").Value = 1679.6666
$ws.Range("K121").Value = 2535
$ws.Range("L121").Value = 5038.9998
$ws.Range("M121").Value = -788
$ws.Range("N121").Value = -8532.9998

$ws.Range("H129").Value = 814.8823
$ws.Range("I129").Value = 441.5
$ws.Range("J129").Value = 929.7692
$ws.Range("K129").Value = 1324.5
$ws.Range("L129").Value = 2789.3076
$ws.Range("M129").Value = 3675.5
$ws.Range("N129").Value = -12789.3076

$ws.Range("H138").Value = 4704.1797
$ws.Range("I138").Value = 934.2917
$ws.Range("J138").Value = 6096.1387
$ws.Range("K138").Value = 2802.8751
$ws.Range("L138").Value = 18288.4161
$ws.Range("M138").Value = 2337.1249
$ws.Range("N138").Value = -28568.4161

$ws.Range("H139").Value = 21999
$ws.Range("J139").Value = 21999
$ws.Range("L139").Value = 21999
$ws.Range("N139").Value = -32279

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("I2").Value = 4300
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 4300
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -4187
$ws.Range("N2").Value = -3226

$ws.Range("H32").Value = 4933.396
$ws.Range("I32").Value = 3450.1707
$ws.Range("J32").Value = 10001.083
$ws.Range("K32").Value = 3450.1707
$ws.Range("L32").Value = 10001.083
$ws.Range("M32").Value = -3163.1707
$ws.Range("N32").Value = -10575.083

$ws.Range("H45").Value = 2107.8462
$ws.Range("I45").Value = 1289.2222
$ws.Range("K45").Value = 1289.2222
$ws.Range("M45").Value = -912.2221999999999

$ws.Range("H97").Value = 1278.2
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 1347.75
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 1347.75
$ws.Range("M97").Value = -504
$ws.Range("N97").Value = -2339.75

$ws.Range("H102").Value = 7409843
$ws.Range("I102").Value = 9261054
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 9261054
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -9259432
$ws.Range("N102").Value = -8244

$ws.Range("I116").Value = 4300
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 4300
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -2006
$ws.Range("N116").Value = -7588

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("I3").Value = 4300
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 4300
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -4186
$ws.Range("N3").Value = -3228

$ws.Range("H20").Value = 1574.6364
$ws.Range("I20").Value = 1291.6
$ws.Range("J20").Value = 2181.1428
$ws.Range("K20").Value = 1291.6
$ws.Range("L20").Value = 2181.1428
$ws.Range("M20").Value = -1044.6
$ws.Range("N20").Value = -2675.1428

$ws.Range("H94").Value = 1633.1515
$ws.Range("I94").Value = 1443.4348
$ws.Range("J94").Value = 2069.5
$ws.Range("K94").Value = 1443.4348
$ws.Range("L94").Value = 2069.5
$ws.Range("M94").Value = -992.4348
$ws.Range("N94").Value = -2971.5

$ws.Range("H124").Value = 52780
$ws.Range("J124").Value = 52780
$ws.Range("L124").Value = 52780
$ws.Range("N124").Value = -62600

$ws.Range("H134").Value = 10537.357
$ws.Range("I134").Value = 29380.75
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 88142.25
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -85607.25
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H16").Value = 792.1429000000001
$ws.Range("I16").Value = 748.75
$ws.Range("J16").Value = 850
$ws.Range("K16").Value = 748.75
$ws.Range("L16").Value = 850
$ws.Range("M16").Value = -461.75
$ws.Range("N16").Value = -1424

$ws.Range("H99").Value = 7498.1665
$ws.Range("I99").Value = 3663
$ws.Range("J99").Value = 11333.333
$ws.Range("K99").Value = 3663
$ws.Range("L99").Value = 11333.333
$ws.Range("M99").Value = -2165
$ws.Range("N99").Value = -14329.333

$ws.Range("H113").Value = 792.1429000000001
$ws.Range("I113").Value = 748.75
$ws.Range("J113").Value = 850
$ws.Range("K113").Value = 748.75
$ws.Range("L113").Value = 850
$ws.Range("M113").Value = 1421.25
$ws.Range("N113").Value = -5190

$ws.Range("H126").Value = 7498.1665
$ws.Range("I126").Value = 3663
$ws.Range("J126").Value = 11333.333
$ws.Range("K126").Value = 10989
$ws.Range("L126").Value = 33999.999
$ws.Range("M126").Value = -8519
$ws.Range("N126").Value = -38939.999

$ws.Range("H132").Value = 2922.375
$ws.Range("I132").Value = 2367.7273
$ws.Range("J132").Value = 4142.6
$ws.Range("K132").Value = 7103.1819
$ws.Range("L132").Value = 12427.8
$ws.Range("M132").Value = -4573.1819
$ws.Range("N132").Value = -17487.8

$ws.Range("H134").Value = 3706.6316
$ws.Range("I134").Value = 3889.7646
$ws.Range("J134").Value = 2150
$ws.Range("K134").Value = 11669.2938
$ws.Range("L134").Value = 6450
$ws.Range("M134").Value = -9134.293799999999
$ws.Range("N134").Value = -11520

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H5").Value = 353842.53
$ws.Range("I5").Value = 480.33334
$ws.Range("K5").Value = 1441.00002
$ws.Range("M5").Value = -1329.00002

$ws.Range("H75").Value = 10209070
$ws.Range("I75").Value = 3006.5
$ws.Range("J75").Value = 14291495
$ws.Range("K75").Value = 9019.5
$ws.Range("L75").Value = 42874485
$ws.Range("M75").Value = -8021.5
$ws.Range("N75").Value = -42876481

$ws.Range("H78").Value = 10209070
$ws.Range("I78").Value = 3006.5
$ws.Range("J78").Value = 14291495
$ws.Range("K78").Value = 27058.5
$ws.Range("L78").Value = 128623455
$ws.Range("M78").Value = -22066.5
$ws.Range("N78").Value = -128633439

$ws.Range("H135").Value = 353842.53
$ws.Range("I135").Value = 480.33334
$ws.Range("K135").Value = 4323.00006
$ws.Range("M135").Value = -1788.00006

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H97").Value = 1003.41174
$ws.Range("I97").Value = 893.13336
$ws.Range("J97").Value = 1830.5
$ws.Range("K97").Value = 893.13336
$ws.Range("L97").Value = 1830.5
$ws.Range("M97").Value = -397.13336
$ws.Range("N97").Value = -2822.5

$ws.Range("H113").Value = 166667840
$ws.Range("I113").Value = 500000640
$ws.Range("J113").Value = 1437.5
$ws.Range("K113").Value = 500000640
$ws.Range("L113").Value = 1437.5
$ws.Range("M113").Value = -499998470
$ws.Range("N113").Value = -5777.5

$ws.Range("H122").Value = 1965965.8
$ws.Range("I122").Value = 2702398.5
$ws.Range("K122").Value = 8107195.5
$ws.Range("M122").Value = -8104745.5

$ws.Range("H123").Value = 10322.75
$ws.Range("J123").Value = 10322.75
$ws.Range("L123").Value = 10322.75
$ws.Range("N123").Value = -15222.75

$ws.Range("H124").Value = 53780
$ws.Range("J124").Value = 53780
$ws.Range("L124").Value = 53780
$ws.Range("N124").Value = -63600

$ws.Range("H132").Value = 7024.7
$ws.Range("I132").Value = 16649.5
$ws.Range("J132").Value = 4618.5
$ws.Range("K132").Value = 49948.5
$ws.Range("L132").Value = 13855.5
$ws.Range("M132").Value = -47418.5
$ws.Range("N132").Value = -18915.5

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H136").Value = 24070.2
$ws.Range("I136").Value = 22550.5
$ws.Range("J136").Value = 26349.75
$ws.Range("K136").Value = 67651.5
$ws.Range("L136").Value = 79049.25
$ws.Range("M136").Value = -65101.5
$ws.Range("N136").Value = -84149.25

$ws.Range("H140").Value = 46171.6
$ws.Range("I140").Value = 35000
$ws.Range("J140").Value = 48964.5
$ws.Range("K140").Value = 35000
$ws.Range("L140").Value = 48964.5
$ws.Range("M140").Value = -29820
$ws.Range("N140").Value = -59324.5

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H122").Value = 2502.862
$ws.Range("I122").Value = 1510.2632
$ws.Range("J122").Value = 4388.8
$ws.Range("K122").Value = 4530.7896
$ws.Range("L122").Value = 13166.4
$ws.Range("M122").Value = -2080.7896
$ws.Range("N122").Value = -18066.4

$ws.Range("H132").Value = 2229.0908
$ws.Range("I132").Value = 1355.375
$ws.Range("K132").Value = 4066.125
$ws.Range("M132").Value = -1536.125

$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280
